# Add a new "hideInContents" column to the survey sheet so notes can be
# hidden from the contents screen, then leave the workbook's view state
# the way Excel would after this edit: cursor parked on the new column's
# next cell, and the settings tab active/selected.

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$survey.Range("F1").Value = "hideInContents"

# Park the selection on the cell just below the new header, matching
# where Excel leaves the cursor after typing a header into a new column.
$survey.Range("F2").Select()

# Switch focus to the settings sheet, making it the active/selected tab.
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
